# edit.ps1
# Applies the "updated pdf with github" change:
#   After the paragraph ending "...they have no meaning on their own."
#   add a new numbered list item:
#       Github: https://github.com/noamzilo/OpenUniversityMaman1
#   (as a real, clickable hyperlink), followed by a blank List Paragraph
#   that takes over the document's auto-managed "_GoBack" edit-position
#   bookmark (Word itself relocates this bookmark whenever new content is
#   typed/inserted - exactly what the target revision shows: the bookmark
#   moves out of the "in this document" run and down into a new, empty
#   paragraph that follows the freshly inserted GitHub link).

$d = $word.ActiveDocument

function Get-ParagraphIndex($rangeStart) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $rangeStart) {
            return $i
        }
    }
    return -1
}

# --- Locate the anchor paragraph robustly (search instead of a hard index) ---
$rng = $d.Content
$found = $rng.Find.Execute("they have no meaning on their own.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the anchor paragraph ('...they have no meaning on their own.')"
}
$rng.Expand(4) | Out-Null   # wdParagraph -> grow the hit to its whole paragraph (incl. mark)
$anchorIndex = Get-ParagraphIndex $rng.Start
if ($anchorIndex -eq -1) {
    throw "Could not resolve the anchor paragraph's index"
}

# --- Insert the new "Github: <link>" list paragraph right after it ---
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter()

$githubIndex = $anchorIndex + 1
$githubPara = $d.Paragraphs.Item($githubIndex)
$githubPara.Range.Text = "Github: https://github.com/noamzilo/OpenUniversityMaman1"

# Convert the URL portion of that text into a real hyperlink, leaving the
# "Github: " label as plain text.
$label = "Github: "
$url = "https://github.com/noamzilo/OpenUniversityMaman1"
$paraStart = $githubPara.Range.Start
$urlRange = $d.Range($paraStart + $label.Length, $paraStart + $label.Length + $url.Length)
$d.Hyperlinks.Add($urlRange, $url) | Out-Null

# --- Insert the trailing empty List Paragraph that will host "_GoBack" ---
$githubPara = $d.Paragraphs.Item($githubIndex)  # re-seat after the hyperlink edit
$githubPara.Range.InsertParagraphAfter()
$blankIndex = $githubIndex + 1
$blankPara = $d.Paragraphs.Item($blankIndex)

# That blank paragraph inherited numbering from the list item above it;
# the target revision shows it as a plain (un-numbered) List Paragraph.
$blankPara.Range.ListFormat.RemoveNumbers()

# Re-home the special, auto-managed "_GoBack" bookmark onto this blank
# paragraph - adding a bookmark with that reserved name moves/replaces
# whichever one already exists elsewhere in the document (Word keeps only
# one "_GoBack" bookmark, tracking the most recent edit location).
$d.Bookmarks.Add("_GoBack", $blankPara.Range) | Out-Null

# --- Register the "Hyperlink" (and "Unresolved Mention") character
#     styles in the style sheet, same as Word does the first time a
#     hyperlink / @mention-aware style is actually used in a document ---
$existingNames = @()
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $existingNames += $d.Styles.Item($i).NameLocal
}

if ($existingNames -notcontains "Hyperlink") {
    $dpf = $d.Styles.Item("DefaultParagraphFont")
    $hl = $d.Styles.Add("Hyperlink", 2)  # wdStyleTypeCharacter
    $hl.BaseStyle = $dpf
    $hl.Priority = 99
    $hl.UnhideWhenUsed = $true
    $hl.Font.Underline = 1               # wdUnderlineSingle
    $hl.Font.Color = 12673797            # BGR-packed 0x0563C1 ("hyperlink" theme blue)
}

if ($existingNames -notcontains "Unresolved Mention") {
    $dpf2 = $d.Styles.Item("DefaultParagraphFont")
    $um = $d.Styles.Add("Unresolved Mention", 2)  # wdStyleTypeCharacter
    $um.BaseStyle = $dpf2
    $um.Priority = 99
    $um.UnhideWhenUsed = $true
    $um.Font.Color = 6316889  # BGR-packed 0x605E5C
}

Write-Output "Github link paragraph inserted; bookmark relocated; styles ensured."
